$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1877.7333
$ws.Cells.Item(19, 10).Value = 2480.2
$ws.Cells.Item(19, 12).Value = 2480.2
$ws.Cells.Item(19, 14).Value = -2830.2
$ws.Cells.Item(55, 8).Value = 165.875
$ws.Cells.Item(55, 10).Value = 201
$ws.Cells.Item(55, 12).Value = 201
$ws.Cells.Item(55, 14).Value = -629
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(106, 8).Value = 2388.111
$ws.Cells.Item(106, 9).Value = 2949
$ws.Cells.Item(106, 10).Value = 1266.3334
$ws.Cells.Item(106, 11).Value = 2949
$ws.Cells.Item(106, 12).Value = 1266.3334
$ws.Cells.Item(106, 13).Value = -2318
$ws.Cells.Item(106, 14).Value = -2528.3334
$ws.Cells.Item(129, 8).Value = 1415.6666
$ws.Cells.Item(129, 10).Value = 1674.25
$ws.Cells.Item(129, 12).Value = 5022.75
$ws.Cells.Item(129, 14).Value = -15022.75
$ws.Cells.Item(135, 8).Value = 740.7646999999999
$ws.Cells.Item(135, 9).Value = 449.44446
$ws.Cells.Item(135, 11).Value = 4045.00014
$ws.Cells.Item(135, 13).Value = -1510.00014
$ws.Cells.Item(138, 8).Value = 2082.0105
$ws.Cells.Item(138, 9).Value = 2325.2285
$ws.Cells.Item(138, 10).Value = 1940.1333
$ws.Cells.Item(138, 11).Value = 6975.685500000001
$ws.Cells.Item(138, 12).Value = 5820.3999
$ws.Cells.Item(138, 13).Value = -1835.685500000001
$ws.Cells.Item(138, 14).Value = -16100.3999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6573.6025
$ws.Cells.Item(32, 9).Value = 4674.294
$ws.Cells.Item(32, 11).Value = 4674.294
$ws.Cells.Item(32, 13).Value = -4387.294
$ws.Cells.Item(45, 8).Value = 1309.2858
$ws.Cells.Item(45, 9).Value = 1028.9
$ws.Cells.Item(45, 10).Value = 1564.1818
$ws.Cells.Item(45, 11).Value = 1028.9
$ws.Cells.Item(45, 12).Value = 1564.1818
$ws.Cells.Item(45, 13).Value = -651.9000000000001
$ws.Cells.Item(45, 14).Value = -2318.1818
$ws.Cells.Item(110, 8).Value = 750
$ws.Cells.Item(110, 9).Value = 750
$ws.Cells.Item(110, 11).Value = 750
$ws.Cells.Item(110, 13).Value = 1295
$ws.Cells.Item(135, 8).Value = 38950
$ws.Cells.Item(135, 10).Value = 38950
$ws.Cells.Item(135, 12).Value = 38950
$ws.Cells.Item(135, 14).Value = -49090

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1314
$ws.Cells.Item(99, 9).Value = 1246.2
$ws.Cells.Item(99, 10).Value = 1398.75
$ws.Cells.Item(99, 11).Value = 1246.2
$ws.Cells.Item(99, 12).Value = 1398.75
$ws.Cells.Item(99, 13).Value = 251.8
$ws.Cells.Item(99, 14).Value = -4394.75
$ws.Cells.Item(100, 8).Value = 30000
$ws.Cells.Item(100, 10).Value = 30000
$ws.Cells.Item(100, 12).Value = 30000
$ws.Cells.Item(100, 14).Value = -32164

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 6215382
$ws.Cells.Item(58, 9).Value = 10873790
$ws.Cells.Item(58, 10).Value = 4171
$ws.Cells.Item(58, 11).Value = 10873790
$ws.Cells.Item(58, 12).Value = 4171
$ws.Cells.Item(58, 13).Value = -10873587
$ws.Cells.Item(58, 14).Value = -4577
$ws.Cells.Item(92, 8).Value = 33747.5
$ws.Cells.Item(92, 10).Value = 33747.5
$ws.Cells.Item(92, 12).Value = 33747.5
$ws.Cells.Item(92, 14).Value = -38739.5
$ws.Cells.Item(95, 8).Value = 22800
$ws.Cells.Item(95, 10).Value = 22800
$ws.Cells.Item(95, 12).Value = 22800
$ws.Cells.Item(95, 14).Value = -28292
$ws.Cells.Item(132, 8).Value = 1322.909
$ws.Cells.Item(132, 9).Value = 952.9474
$ws.Cells.Item(132, 10).Value = 3666
$ws.Cells.Item(132, 11).Value = 2858.8422
$ws.Cells.Item(132, 12).Value = 10998
$ws.Cells.Item(132, 13).Value = -328.8422
$ws.Cells.Item(132, 14).Value = -16058
$ws.Cells.Item(136, 8).Value = 6215382
$ws.Cells.Item(136, 9).Value = 10873790
$ws.Cells.Item(136, 10).Value = 4171
$ws.Cells.Item(136, 11).Value = 32621370
$ws.Cells.Item(136, 12).Value = 12513
$ws.Cells.Item(136, 13).Value = -32618820
$ws.Cells.Item(136, 14).Value = -17613

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2593571
$ws.Cells.Item(4, 9).Value = 4000119.8
$ws.Cells.Item(4, 11).Value = 12000359.4
$ws.Cells.Item(4, 13).Value = -12000247.4
$ws.Cells.Item(23, 8).Value = 909271.2
$ws.Cells.Item(23, 10).Value = 1666945.9
$ws.Cells.Item(23, 12).Value = 5000837.699999999
$ws.Cells.Item(23, 14).Value = -5001307.699999999
$ws.Cells.Item(107, 8).Value = 453.6154
$ws.Cells.Item(107, 10).Value = 427
$ws.Cells.Item(107, 12).Value = 1281
$ws.Cells.Item(107, 14).Value = -5121
$ws.Cells.Item(131, 8).Value = 12673.717
$ws.Cells.Item(131, 10).Value = 13047.062
$ws.Cells.Item(131, 12).Value = 39141.186
$ws.Cells.Item(131, 14).Value = -49221.186
$ws.Cells.Item(136, 8).Value = 1291.8462
$ws.Cells.Item(136, 9).Value = 1291.8462
$ws.Cells.Item(136, 11).Value = 3875.5386
$ws.Cells.Item(136, 13).Value = 1224.4614
$ws.Cells.Item(138, 8).Value = 2276.7693
$ws.Cells.Item(138, 9).Value = 1586.2727
$ws.Cells.Item(138, 10).Value = 6074.5
$ws.Cells.Item(138, 11).Value = 4758.8181
$ws.Cells.Item(138, 12).Value = 18223.5
$ws.Cells.Item(138, 13).Value = 381.1818999999996
$ws.Cells.Item(138, 14).Value = -28503.5
$ws.Cells.Item(139, 8).Value = 6808.65
$ws.Cells.Item(139, 9).Value = 7061.737
$ws.Cells.Item(139, 11).Value = 21185.211
$ws.Cells.Item(139, 13).Value = -16045.211
$ws.Cells.Item(140, 8).Value = 2322.7646
$ws.Cells.Item(140, 9).Value = 1435.0714
$ws.Cells.Item(140, 11).Value = 4305.2142
$ws.Cells.Item(140, 13).Value = 874.7857999999997
$ws.Cells.Item(141, 8).Value = 3958.5454
$ws.Cells.Item(141, 9).Value = 3184.75
$ws.Cells.Item(141, 11).Value = 9554.25
$ws.Cells.Item(141, 13).Value = -4374.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1215.2
$ws.Cells.Item(113, 9).Value = 1037.4286
$ws.Cells.Item(113, 10).Value = 1370.75
$ws.Cells.Item(113, 11).Value = 1037.4286
$ws.Cells.Item(113, 12).Value = 1370.75
$ws.Cells.Item(113, 13).Value = 1132.5714
$ws.Cells.Item(113, 14).Value = -5710.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 15775.3
$ws.Cells.Item(40, 9).Value = 21382.334
$ws.Cells.Item(40, 11).Value = 21382.334
$ws.Cells.Item(40, 13).Value = -21246.334
$ws.Cells.Item(43, 8).Value = 14405.6
$ws.Cells.Item(43, 10).Value = 14405.6
$ws.Cells.Item(43, 12).Value = 14405.6
$ws.Cells.Item(43, 14).Value = -14791.6
$ws.Cells.Item(61, 8).Value = 4506
$ws.Cells.Item(61, 9).Value = 4210.5
$ws.Cells.Item(61, 10).Value = 4900
$ws.Cells.Item(61, 11).Value = 4210.5
$ws.Cells.Item(61, 12).Value = 4900
$ws.Cells.Item(61, 13).Value = -4008.5
$ws.Cells.Item(61, 14).Value = -5304
$ws.Cells.Item(93, 8).Value = 23810932
$ws.Cells.Item(93, 9).Value = 783
$ws.Cells.Item(93, 10).Value = 66669200
$ws.Cells.Item(93, 11).Value = 783
$ws.Cells.Item(93, 12).Value = 66669200
$ws.Cells.Item(93, 13).Value = 465
$ws.Cells.Item(93, 14).Value = -66671696
$ws.Cells.Item(100, 8).Value = 1523
$ws.Cells.Item(100, 9).Value = 1483.4286
$ws.Cells.Item(100, 11).Value = 1483.4286
$ws.Cells.Item(100, 13).Value = -942.4286
$ws.Cells.Item(113, 8).Value = 4506
$ws.Cells.Item(113, 9).Value = 4210.5
$ws.Cells.Item(113, 10).Value = 4900
$ws.Cells.Item(113, 11).Value = 4210.5
$ws.Cells.Item(113, 12).Value = 4900
$ws.Cells.Item(113, 13).Value = -2040.5
$ws.Cells.Item(113, 14).Value = -9240
$ws.Cells.Item(132, 8).Value = 1984.2826
$ws.Cells.Item(132, 9).Value = 1482.8182
$ws.Cells.Item(132, 11).Value = 4448.4546
$ws.Cells.Item(132, 13).Value = -1918.4546

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1213.6666
$ws.Cells.Item(100, 9).Value = 1009
$ws.Cells.Item(100, 11).Value = 2018
$ws.Cells.Item(100, 13).Value = -1477

Write-Output "Applied all cell updates"